# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund holdings detail) positioned right
#    after "2021-Q4" and before "总计".
# 2) Insert a new first data row into "总计" summarising the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet right after "2021-Q4"
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Pick up the same formatting (bold / centered / bordered) used by the
# existing quarterly sheets' header row and index column.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: index, code, name, scale, stock position, position pct, market value, rank
$rows = @(
    @(0, "162006", "长城久富核心成长混合(LOF)", "19.40", "79.43", "4.83", "0.9370", 3),
    @(1, "010846", "南方卓越优选3个月持有期混合A", "26.01", "60.36", "2.92", "0.7595", 5),
    @(2, "013049", "兴业能源革新股票A", "8.80", "90.46", "3.21", "0.2825", 9),
    @(3, "013050", "兴业能源革新股票C", "3.84", "90.46", "3.21", "0.1233", 9),
    @(4, "010847", "南方卓越优选3个月持有期混合C", "3.42", "60.36", "2.92", "0.0999", 5),
    @(5, "000976", "长城新兴产业灵活配置混合", "1.02", "77.95", "4.72", "0.0481", 3),
    @(6, "004703", "南方兴盛先锋灵活配置混合", "1.09", "53.08", "2.27", "0.0247", 10)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    # Codes / name / ratio columns are stored as text in the source data
    # (leading zeros in fund codes must be preserved), so force text format
    # before writing.
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[1]

    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[3]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[4]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[5]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new first data row for 2022-Q1
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The row-insert leaves a stray inherited text format on the shifted-in
# blank row; clear it so the new data cells match the plain (unstyled) data
# cells used elsewhere in this sheet.
$total.Range("B2:D2").ClearFormats()

# Re-apply the bold/centered/bordered index-column style (used by every
# other row in column A) to the new row's A cell.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 2.28

# The index column (A) is a simple positional counter; renumber the rows
# that were pushed down by the insert.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
